# Updates the cryptos list data (prices / 1h volume % / and two swapped
# coin rows) in the active worksheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '30.649.84'
    'E2' = '  +0.81%  '
    'D3' = '2.089.00'
    'E3' = '  +4.59%  '
    'E4' = '  +0.18%  '
    'D5' = '329.92'
    'E5' = '  +1.83%  '
    'D6' = '1.006'
    'E6' = '  +0.51%  '
    'D7' = '0.5233'
    'E7' = '  +2.55%  '
    'D8' = '0.4302'
    'E8' = '  +4.23%  '
    'D9' = '0.08829'
    'E9' = '  +1.78%  '
    'D10' = '46.33'
    'E10' = '  +8.28%  '
    'D11' = '1.160'
    'E11' = '  +2.60%  '
    'D12' = '24.32'
    'E12' = '  -1.35%  '
    'D13' = '2.088.70'
    'E13' = '  +4.47%  '
    'D14' = '6.653'
    'E14' = '  +1.64%  '
    'D15' = '7.776'
    'E15' = '  +4.67%  '
    'D16' = '96.39'
    'E16' = '  +2.54%  '
    'D17' = '1.009'
    'E17' = '  +0.50%  '
    'E18' = '  +0.83%  '
    'D19' = '0.06645'
    'E19' = '  +2.23%  '
    'D20' = '18.82'
    'E20' = '  -0.41%  '
    'D21' = '1.006'
    'E21' = '  +0.71%  '
    'D22' = '6.270'
    'E22' = '  +1.97%  '
    'D23' = '30.625.99'
    'E23' = '  +0.36%  '
    'D24' = '12.35'
    'E24' = '  +4.78%  '
    'B25' = 'Toncoin'
    'C25' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D25' = '2.297'
    'E25' = '  +3.84%  '
    'B26' = 'WrappedliquidstakedEther2.0'
    'C26' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D26' = '2.344.74'
    'E26' = '  +5.06%  '
    'D27' = '22.25'
    'E27' = '  -0.43%  '
    'D28' = '2.525'
    'E28' = '  +4.91%  '
    'D29' = '161.95'
    'E29' = '  -0.59%  '
    'D30' = '132.30'
    'E30' = '  +0.48%  '
    'D31' = '1.191'
    'E31' = '  +4.89%  '
    'D32' = '0.1073'
    'E32' = '  +1.90%  '
    'D33' = '6.102'
    'E33' = '  +0.51%  '
    'D34' = '1.533'
    'E34' = '  +15.45%  '
    'D35' = '3.838'
    'E35' = '  -0.25%  '
    'D36' = '0.02570'
    'E36' = '  +2.23%  '
    'D37' = '9.638'
    'E37' = '  +7.37%  '
    'D38' = '5.469'
    'E38' = '  +2.37%  '
    'D39' = '0.06672'
    'E39' = '  +1.20%  '
    'D40' = '0.2252'
    'D41' = '12.48'
    'E41' = '  +2.07%  '
    'D42' = '0.6708'
    'E42' = '  +1.59%  '
    'D43' = '1.241'
    'E43' = '  +0.52%  '
    'D44' = '1.004'
    'E44' = '  +0.42%  '
    'D45' = '0.6333'
    'E45' = '  +3.08%  '
    'D46' = '13.86'
    'E46' = '  +1.33%  '
    'D47' = '2.194'
    'E47' = '  -0.21%  '
    'D48' = '3.621'
    'E48' = '  -0.99%  '
    'D49' = '1.234'
    'E49' = '  -1.94%  '
    'B50' = 'Aave'
    'C50' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D50' = '82.37'
    'E50' = '  +2.81%  '
    'B51' = 'WEMIXTOKEN'
    'C51' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D51' = '1.186'
    'E51' = '  +6.61%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $val = $updates[$addr]
    # Preserve text storage for numeric-looking values (prices) by
    # forcing them in as text, same as the original inline-string data,
    # instead of letting Excel auto-convert them to numbers.
    if ($addr -match '^D\d+$' -and $val -match '^-?[0-9.]+$') {
        $cell.Value = "'" + $val
    } else {
        $cell.Value = $val
    }
}
